$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.409.27"
$ws.Range("E2").Value = "  +1.41%  "
$ws.Range("D3").Value = "1.824.99"
$ws.Range("E3").Value = "  +2.44%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.41"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("E7").Value = "  -0.97%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4028"
$ws.Range("E8").Value = "  +6.77%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07614"
$ws.Range("E9").Value = "  +2.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.83"
$ws.Range("E10").Value = "  +0.36%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.109"
$ws.Range("E11").Value = "  +1.35%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.336"
$ws.Range("E12").Value = "  +4.21%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.001"
$ws.Range("E13").Value = "  +0.18%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.621"
$ws.Range("E14").Value = "  +5.53%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.91"
$ws.Range("E15").Value = "  +2.08%  "
$ws.Range("D16").Value = "1.828.47"
$ws.Range("E16").Value = "  +3.00%  "
$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "89.52"
$ws.Range("E17").Value = "  +1.30%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001076"
$ws.Range("E18").Value = "  +2.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06594"
$ws.Range("E19").Value = "  +2.44%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.68"
$ws.Range("E20").Value = "  +2.55%  "
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.082"
$ws.Range("E22").Value = "  +3.47%  "
$ws.Range("D23").Value = "28.413.13"
$ws.Range("E23").Value = "  +1.39%  "
$ws.Range("E24").Value = "  +0.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.217"
$ws.Range("E25").Value = "  +6.17%  "
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.455"
$ws.Range("E26").Value = "  +6.98%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "157.67"
$ws.Range("E27").Value = "  +0.99%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.63"
$ws.Range("E28").Value = "  +2.01%  "
$ws.Range("D29").Value = "2.039.46"
$ws.Range("E29").Value = "  +3.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "123.99"
$ws.Range("E30").Value = "  +3.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.122"
$ws.Range("E31").Value = "  +1.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1103"
$ws.Range("E32").Value = "  +4.84%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.660"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07462"
$ws.Range("E34").Value = "  +16.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.645"
$ws.Range("E35").Value = "  +0.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2231"
$ws.Range("E36").Value = "  -1.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02341"
$ws.Range("E37").Value = "  +2.96%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.912"
$ws.Range("E38").Value = "  +5.77%  "
$ws.Range("E39").Value = "  +4.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.31"
$ws.Range("E40").Value = "  +1.76%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6255"
$ws.Range("E41").Value = "  +1.80%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.177"
$ws.Range("E42").Value = "  +0.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9996"
$ws.Range("E43").Value = "  +0.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.393"
$ws.Range("E44").Value = "  -2.51%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.45"
$ws.Range("E45").Value = "  +1.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.704"
$ws.Range("E46").Value = "  +1.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5842"
$ws.Range("E47").Value = "  +1.60%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.86"
$ws.Range("E48").Value = "  -1.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.992"
$ws.Range("E49").Value = "  +3.46%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.201"
$ws.Range("E50").Value = "  +1.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06890"
$ws.Range("E51").Value = "  +1.48%  "
